$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of the "CDO Override reason" header to "CDO Override Reason"
$ws.Range("AF1").Value = "CDO Override Reason"

# Move/refresh the active selection to AF1 (the cell that was just edited)
$ws.Range("AF1").Select()
